$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Do you remember when you were a kid and you saw that Zinedine Zidane had moved to Real Madrid for £45m? It seemed like so much money at the time - but it''s fairly standard these days for chairmen of big clubs across Europe to PayPal each other figures that large for players. £45m these days buys you the right leg of some players. Tens of millions of pounds is still a lot of money in the real world - and yet we baulk at the idea of receiving £20m for a player worth 30, £50m for a player worth 100.'
$ws.Range("C2").Value = 'Football'
$ws.Range("D2").Value = 100

$ws.Range("B3").Value = 'Mohamed Salah''s sensational second-half goal divided the sides at Anfield, but there was division in the studio after the full-time whistle as Carragher praised Liverpool for stifling Kevin De Bruyne but Neville disagreed and blamed City''s system. Guardiola elected to start with three centre-backs for the first time in the Premier League this season with De Bruyne stationed alongside Ilkay Gundogan behind Erling Haaland. In term of work rate, Liverpool clocked 118.59 km - a distance only surpassed once by Klopp''s side over the past three league campaigns.'
$ws.Range("C3").Value = 'Football'
$ws.Range("D3").Value = 94

$ws.Range("B4").Value = 'Manchester United were faster, sharper, smarter, and classier than Nottingham Forest in this seventh straight win over their visitors. This was approaching the complete performance: a blend of measured attacking and rapid breaks and the constant harrying of an opponent who ended exhausted and demoralised. Given United had scored only 20 Premier League goals at kickoff and there is limited finance in January to replace Cristiano Ronaldo, a strike apiece from Marcus Rashford and Anthony Martial before Fred’s late clincher was as fine a tonic as the three points which keep them on the heels of Tottenham.'
$ws.Range("C4").Value = 'Football'
$ws.Range("D4").Value = 117

$ws.Range("B5").Value = 'Andy Murray and his brother Jamie beat Neal Skupski and Dan Evans in straight sets in Aberdeen but Scotland fall to 8-4 defeat against England; Andy Murray will head to the Australian Open encouraged by his performances. Andy and Jamie Murray were victorious in possibly their final match together in Scotland but England walked away with the inaugural Schroders Battle of the Brits trophy. Organised by Jamie Murray, the Scotland versus England clash proved a hit in Aberdeen, with around 20,000 fans attending the P&J Live arena across two days.'
$ws.Range("C5").Value = 'Tenis'
$ws.Range("D5").Value = 93

$ws.Range("B6").Value = 'The 6ft 1in right-handed Austrian has always had the game to do well on the surface with a booming backhand, stamina, graceful movement and aggressive play from the baseline. To prove a point, he conquered Nadal in the Barcelona Open semi-finals to become the first player other than Novak Djokovic to defeat the Spaniard on clay four times before capping an impressive week by crushing Daniil Medvedev 6-4 6-0 to win the title. He takes on 33-year-old Nadal in a repeat of last year''s Roland Garros final but the Austrian will be feeling far more confident he can give the Spaniard a run for his money after ending Novak Djokovic''s 26-match winning streak at Grand Slams. He is Austria''s biggest tennis star since Thomas Muster, and now the 25-year-old is playing a free-flowing game under the tutelage of Olympic gold medallist Nicolas Massu.'
$ws.Range("C6").Value = 'Tenis'
$ws.Range("D6").Value = 156

$ws.Range("B7").Value = 'Whether Tiger Woods makes it to Italy or not for next year''s Ryder Cup, he will be an integral member of the U.S. team, captain Zach Johnson said Tuesday at the year-to-go ceremonies. Woods was a player on eight Ryder Cup teams and a vice-captain in 2018, before breaking bones in his right leg and ankle in a February 2021 car crash outside Los Angeles. The U.S. romped to a record 19-9 rout of Europe at Whistling Straits, Wisconsin, last year, and also won the Presidents Cup comfortably last month.'
$ws.Range("C7").Value = 'Golf'
$ws.Range("D7").Value = 97

$ws.Range("B8").Value = 'Adam Scott has confirmed he will join compatriots Cameron Smith and Marc Leishman at the Australian PGA Championship in November; LIV Golfers will be allowed to compete despite the event being co-sanctioned by the DP World Tour. Former Masters champion Adam Scott said he does not see the controversial LIV Series as ''pure evil'' for the game of golf and called on the Saudi-backed circuit and the PGA Tour to move on from their feud. LIV Golf has lured away some of golf''s biggest names with huge sums of money, while those who joined the rebel circuit were suspended by the PGA Tour. The breakaway circuit has filed a lawsuit accusing the PGA Tour of antitrust violations, while the PGA Tour Tour has filed a counterclaim. '
$ws.Range("C8").Value = 'Golf'
$ws.Range("D8").Value = 129

$ws.Range("B9").Value = 'Well, for me it''s Devin Booker, being a Suns draftee, it''s been so fun watching him develop from a bench player to one of the best players in the league. I have always loved high-scoring tw-guards and Booker fits that mould perfectly. Though he comes across cocky at times, he has that superstar confidence that trickles down through to the rest of the team and this confidence is what has ushered much of the Suns'' recent success.'
$ws.Range("C9").Value = 'Basketball'
$ws.Range("D9").Value = 81

$ws.Range("B10").Value = 'Coming off their NBA Finals run in 2021, the Suns were expected to potentially get back there again in 2022, but ultimately flamed out against Luka Doncic and the Dallas Mavericks in the second round. It was a painful loss for, not only us Suns fans but all NBA fans, as a Suns vs Warriors series in the Western Conference Finals would have been immense. Although the season did end on that sour note, it was in fact a historically great season for Phoenix. They were far and away the most consistent team all season boasting a top-five offensive and defensive rating. This efficiency helped them win a franchise-record 64 games, leading the NBA, with the next closest team being the Grizzlies with 56 wins.'
$ws.Range("C10").Value = 'Basketball'
$ws.Range("D10").Value = 127

$ws.Range("B11").Value = 'James Harden has vowed to get back to his high-scoring best for the Philadelphia 76ers after instructing president of basketball operations Daryl Morey to build a contender around him and Joel Embiid. The former league MVP said he would be happy to take a pay cut to allow the team to build a title-contending roster and he''s been true to his word, agreeing a cut-price deal which has helped allow Morey time to put the pieces together to allow the Sixers to go full-tilt at winning a championship. Harden has now fully agreed to a two-year, $68.6m deal with the Philadelphia 76ers, ESPN reported on Wednesday. He averaged 22 points this past season for Brooklyn and Philadelphia, the lowest since he became a starter in the 2012-13 season. He turns 33 in August.'
$ws.Range("C11").Value = 'Basketball'
$ws.Range("D11").Value = 141

# D11 needs centered alignment to match style used by the rest of column D
$ws.Range("D11").HorizontalAlignment = -4108
$ws.Range("D11").VerticalAlignment = -4108

# Update active selection from G12 to F12
[void]$ws.Range("F12").Select()
